# Applies the data refresh described in the diff:
#  - Row 2 values change (B2:F2)
#  - Rows 3-8: column B cleared, column C values replaced
#  - Row 3: D3 cleared; Rows 3-4: F cleared
#  - Rows 9-15: column C values replaced only
#  - New rows 16-19 appended with index numbers 14-17 and new tickers in column C

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2, 2).Value = "NSE:BSE"
$ws.Cells.Item(2, 3).Value = "NSE:3MINDIA"
$ws.Cells.Item(2, 4).Value = "NSE:IDEA"
$ws.Cells.Item(2, 5).Value = "NSE:COLPAL"
$ws.Cells.Item(2, 6).Value = "NSE:BSE"

# --- Row 3 ---
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = "NSE:63MOONS"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 6).Value = ""

# --- Row 4 ---
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = "NSE:ASALCBR"
$ws.Cells.Item(4, 6).Value = ""

# --- Row 5 ---
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 3).Value = "NSE:BHAGCHEM"

# --- Row 6 ---
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(6, 3).Value = "NSE:CCL"

# --- Row 7 ---
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(7, 3).Value = "NSE:CHEMFAB"

# --- Row 8 ---
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 3).Value = "NSE:CMSINFO"

# --- Row 9 ---
$ws.Cells.Item(9, 3).Value = "NSE:DNAMEDIA"

# --- Row 10 ---
$ws.Cells.Item(10, 3).Value = "NSE:EMMBI"

# --- Row 11 ---
$ws.Cells.Item(11, 3).Value = "NSE:FINOPB"

# --- Row 12 ---
$ws.Cells.Item(12, 3).Value = "NSE:HDFCGROWTH"

# --- Row 13 ---
$ws.Cells.Item(13, 3).Value = "NSE:INDOCO"

# --- Row 14 ---
$ws.Cells.Item(14, 3).Value = "NSE:JINDWORLD"

# --- Row 15 ---
$ws.Cells.Item(15, 3).Value = "NSE:MARATHON"

# --- New rows 16-19 ---
# Copy formatting from row 15 down into the new rows first so that column A
# keeps the bold/centered/bordered style (s="1") used throughout the table.
$ws.Range("A15:F15").Copy($ws.Range("A16:F16"))
$ws.Range("A15:F15").Copy($ws.Range("A17:F17"))
$ws.Range("A15:F15").Copy($ws.Range("A18:F18"))
$ws.Range("A15:F15").Copy($ws.Range("A19:F19"))

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = "NSE:MARICO"

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = "NSE:POWERMECH"

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 3).Value = "NSE:RAMKY"

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 3).Value = "NSE:ROTO"
